# Update the build-version timestamp embedded in several cells across the
# "About" and "Boundaries and methane sources" worksheets.
#
# Old: "January 30 2026 16.19.47 EST"
# New: "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet -------------------------------------------------------
$a2 = $wsAbout.Range("A2")
$a2Text = [string]$a2.Value()
$a2.Value = $a2Text.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6")
$a6Text = [string]$a6.Value()
$a6.Value = $a6Text.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet -------------------------------
for ($r = 2; $r -le 11; $r++) {
    $cell = $wsData.Range("S$r")
    $cellText = [string]$cell.Value()
    $cell.Value = $cellText.Replace($oldStamp, $newStamp)
}
